$wb = $excel.ActiveWorkbook

# --- OFF sheet: update row 2 (H) ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 545
$wsOff.Range("C2").Value = 416
$wsOff.Range("D2").Value = 151
$wsOff.Range("E2").Value = 74

# --- DEF sheet: update row 2 (H) ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 574
$wsDef.Range("C2").Value = 404
$wsDef.Range("D2").Value = 117
$wsDef.Range("E2").Value = 48
$wsDef.Range("F2").Value = 9
